$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 9 entirely (company list shrinks from 8 data rows to 7)
$ws.Rows.Item(9).Delete()

# --- Row 2 ---
$row2 = New-Object 'object[,]' 1,43
$row2[0,0] = "United Kingdom"
$row2[0,1] = "'6"
$row2[0,2] = "Insurance (Prop/Cas.)"
$row2[0,3] = 0.0789
$row2[0,4] = 0.0354
$row2[0,5] = 0.03154
$row2[0,6] = 0.165263497085274
$row2[0,7] = 0.165263497085274
$row2[0,8] = 0.1412133462225738
$row2[0,9] = 0.117746174684161
$row2[0,10] = 1676.1
$row2[0,11] = 0.09819843571491343
$row2[0,12] = 687.746
$row2[0,13] = 0.02223822442387224
$row2[0,14] = 0.4103251595966828
$row2[0,15] = 598.65
$row2[0,16] = 0.01935731076785778
$row2[0,17] = 0.3571684267048506
$row2[0,18] = 89.096
$row2[0,19] = 0.1295478272501767
$row2[0,20] = 3907.6
$row2[0,21] = 0.1263520046044952
$row2[0,22] = 0.1499925335689465
$row2[0,23] = 0.04679304211923487
$row2[0,24] = 0.1031994914497116
$row2[0,25] = 1.601682561971419
$row2[0,26] = 0.3548463041953803
$row2[0,27] = 0.04445638117149156
$row2[0,28] = 0.3105436173092129
$row2[0,29] = 3558.338
$row2[0,30] = 0
$row2[0,31] = 3558.338
$row2[0,32] = -349.2620000000002
$row2[0,33] = 0.1031861781469186
$row2[0,34] = 0.2133406414648047
$row2[0,35] = -0.01142236210060635
$row2[0,36] = -0.02734690215070418
$row2[0,37] = 144.363
$row2[0,38] = 144.363
$row2[0,39] = 1.39778371371332
$row2[0,40] = 16.69610634303804
$row2[0,41] = -0.1371968417331187
$row2[0,42] = 16.69610634303804
$ws.Range("A2:AQ2").Value = $row2

# --- Row 3 ---
$row3 = New-Object 'object[,]' 1,43
$row3[0,0] = "United Kingdom"
$row3[0,1] = "Personal Group Holdings Plc (AIM:PGH)"
$row3[0,2] = "Insurance (Prop/Cas.)"
$row3[0,3] = 0.0794
$row3[0,4] = 0.0354
$row3[0,5] = $null
$row3[0,6] = 0.1528878822197056
$row3[0,7] = 0.1528878822197056
$row3[0,8] = 0.1415628539071348
$row3[0,9] = 0.1167083070379432
$row3[0,10] = 10.8
$row3[0,11] = 0.1223103057757644
$row3[0,12] = 7.419
$row3[0,13] = 0.08469178082191781
$row3[0,14] = 0.6869444444444444
$row3[0,15] = 7.35
$row3[0,16] = 0.0839041095890411
$row3[0,17] = 0.6805555555555555
$row3[0,18] = 0.06899999999999995
$row3[0,19] = 0.009300444803881919
$row3[0,20] = 19.7
$row3[0,21] = 0.2248858447488584
$row3[0,22] = 0.2337662337662338
$row3[0,23] = 0.0440802498927258
$row3[0,24] = 0.189685983873508
$row3[0,25] = 8.908393866020981
$row3[0,26] = 1.039683566530506
$row3[0,27] = 0.04404108254468058
$row3[0,28] = 0.9956424839858259
$row3[0,29] = 0.142
$row3[0,30] = 0
$row3[0,31] = 0.142
$row3[0,32] = -19.558
$row3[0,33] = 0.001618381162955027
$row3[0,34] = 0.002919287858229513
$row3[0,35] = -0.2874401105199729
$row3[0,36] = -0.6757653237509502
$row3[0,37] = 0.132
$row3[0,38] = 0.132
$row3[0,39] = 0.01092307692307692
$row3[0,40] = 94.69696969696969
$row3[0,41] = -1.504461538461539
$row3[0,42] = 94.69696969696969
$ws.Range("A3:AQ3").Value = $row3

# --- Row 4 ---
$row4 = New-Object 'object[,]' 1,43
$row4[0,0] = "United Kingdom"
$row4[0,1] = "Admiral Group plc (LSE:ADM)"
$row4[0,2] = "Insurance (Prop/Cas.)"
$row4[0,3] = 0.0956
$row4[0,4] = 0.11
$row4[0,5] = -0.00622
$row4[0,6] = 0.4224873213462425
$row4[0,7] = 0.4224873213462425
$row4[0,8] = 0.4355117565698479
$row4[0,9] = 0.3615265833158811
$row4[0,10] = 610.3
$row4[0,11] = 0.3517173812816966
$row4[0,12] = 236.8
$row4[0,13] = 0.02073355456129445
$row4[0,14] = 0.3880058987383255
$row4[0,15] = 236.8
$row4[0,16] = 0.02073355456129445
$row4[0,17] = 0.3880058987383255
$row4[0,18] = 0
$row4[0,19] = 0
$row4[0,20] = 485.2
$row4[0,21] = 0.04248277311292257
$row4[0,22] = 0.5921218589308237
$row4[0,23] = 0.04640396814756931
$row4[0,24] = 0.5457178907832544
$row4[0,25] = 1.485362095531587
$row4[0,26] = 0.5369978833844522
$row4[0,27] = 0.04434370664451491
$row4[0,28] = 0.4926541767399373
$row4[0,29] = 961.8
$row4[0,30] = 0
$row4[0,31] = 961.8
$row4[0,32] = 476.6
$row4[0,33] = 0.07767162780931769
$row4[0,34] = 0.4250110472823685
$row4[0,35] = 0.04005816250199618
$row4[0,36] = 0.2680841489481381
$row4[0,37] = 26.7
$row4[0,38] = 26.7
$row4[0,39] = 1.245854922279793
$row4[0,40] = 28.30337078651686
$row4[0,41] = 0.6173575129533678
$row4[0,42] = 28.30337078651686
$ws.Range("A4:AQ4").Value = $row4

# --- Row 5 ---
$row5 = New-Object 'object[,]' 1,43
$row5[0,0] = "United Kingdom"
$row5[0,1] = "Sabre Insurance Group plc (LSE:SBRE)"
$row5[0,2] = "Insurance (Prop/Cas.)"
$row5[0,3] = $null
$row5[0,4] = $null
$row5[0,5] = $null
$row5[0,6] = 0.3302509907529723
$row5[0,7] = 0.3302509907529723
$row5[0,8] = 0.2932628797886394
$row5[0,9] = 0.23778071334214
$row5[0,10] = 53.9
$row5[0,11] = 0.2373403786878027
$row5[0,12] = 21.627
$row5[0,13] = 0.02300255264837268
$row5[0,14] = 0.4012430426716141
$row5[0,15] = 21
$row5[0,16] = 0.0223356732610083
$row5[0,17] = 0.3896103896103896
$row5[0,18] = 0.6269999999999989
$row5[0,19] = 0.02899153835483418
$row5[0,20] = 46.1
$row5[0,21] = 0.04903212082535631
$row5[0,22] = 0.164279183175861
$row5[0,23] = 0.04404649295215926
$row5[0,24] = 0.1202326902237018
$row5[0,25] = 2.138538900502852
$row5[0,26] = 0.508503305271484
$row5[0,27] = 0.04403631807386676
$row5[0,28] = 0.4644669871976173
$row5[0,29] = 0.396
$row5[0,30] = 0
$row5[0,31] = 0.396
$row5[0,32] = -45.704
$row5[0,33] = 0.0004210096577063905
$row5[0,34] = 0.001175436930091185
$row5[0,35] = -0.05109469466604658
$row5[0,36] = -0.1571685992929752
$row5[0,37] = 0.031
$row5[0,38] = 0.031
$row5[0,39] = 0.005928143712574851
$row5[0,40] = 2148.387096774193
$row5[0,41] = -0.6841916167664671
$row5[0,42] = 2148.387096774193
$ws.Range("A5:AQ5").Value = $row5

# --- Row 6 ---
$row6 = New-Object 'object[,]' 1,43
$row6[0,0] = "United Kingdom"
$row6[0,1] = "Direct Line Insurance Group plc (LSE:DLG)"
$row6[0,2] = "Insurance (Prop/Cas.)"
$row6[0,3] = 0.00157
$row6[0,4] = -0.08500000000000001
$row6[0,5] = -0.0137
$row6[0,6] = 0.2022021526660893
$row6[0,7] = 0.2022021526660893
$row6[0,8] = 0.16966472844241
$row6[0,9] = 0.1402426091173144
$row6[0,10] = 496.1
$row6[0,11] = 0.1227514536681925
$row6[0,12] = 210.5
$row6[0,13] = 0.03571428571428571
$row6[0,14] = 0.4243096149969764
$row6[0,15] = 142.6
$row6[0,16] = 0.02419409569053274
$row6[0,17] = 0.2874420479741987
$row6[0,18] = 67.90000000000001
$row6[0,19] = 0.3225653206650831
$row6[0,20] = 1857.3
$row6[0,21] = 0.3151170682049542
$row6[0,22] = 0.1357058839620319
$row6[0,23] = 0.04839953925056728
$row6[0,24] = 0.08730634471146462
$row6[0,25] = 1.434580434473946
$row6[0,26] = 0.2011893031192767
$row6[0,27] = 0.04456905569846822
$row6[0,28] = 0.1566202474208084
$row6[0,29] = 914.4
$row6[0,30] = 0
$row6[0,31] = 914.4
$row6[0,32] = -942.9
$row6[0,33] = 0.1343046824510898
$row6[0,34] = 0.1922421948912015
$row6[0,35] = -0.1904425279230878
$row6[0,36] = -0.3252276490066225
$row6[0,37] = 32.8
$row6[0,38] = 32.8
$row6[0,39] = 1.207766477347774
$row6[0,40] = 20.90548780487805
$row6[0,41] = -1.245410117553824
$row6[0,42] = 20.90548780487805
$ws.Range("A6:AQ6").Value = $row6

# --- Row 7 ---
$row7 = New-Object 'object[,]' 1,43
$row7[0,0] = "United Kingdom"
$row7[0,1] = "RSA Insurance Group plc (LSE:RSA)"
$row7[0,2] = "Insurance (Prop/Cas.)"
$row7[0,3] = -0.00357
$row7[0,4] = 0.0425
$row7[0,5] = 0.0693
$row7[0,6] = 0.1054720862253552
$row7[0,7] = 0.1054720862253552
$row7[0,8] = 0.09219184179187065
$row7[0,9] = 0.07049688343331541
$row7[0,10] = 422.2
$row7[0,11] = 0.0507873115925467
$row7[0,12] = 113.9
$row7[0,13] = 0.01189730091083814
$row7[0,14] = 0.2697773567029844
$row7[0,15] = 113.9
$row7[0,16] = 0.01189730091083814
$row7[0,17] = 0.2697773567029844
$row7[0,18] = 0
$row7[0,19] = 0
$row7[0,20] = 1149
$row7[0,21] = 0.1200175482577087
$row7[0,22] = 0.08217684956303405
$row7[0,23] = 0.04718211609090042
$row7[0,24] = 0.03499473347213362
$row7[0,25] = 1.702283198525648
$row7[0,26] = 0.1200056602169539
$row7[0,27] = 0.04485063647243124
$row7[0,28] = 0.07515502374452265
$row7[0,29] = 1071
$row7[0,30] = 0
$row7[0,31] = 1071
$row7[0,32] = -78
$row7[0,33] = 0.1006143960317908
$row7[0,34] = 0.1566728594625433
$row7[0,35] = -0.008214330847971692
$row7[0,36] = -0.01371573264871899
$row7[0,37] = 49.5
$row7[0,38] = 49.5
$row7[0,39] = 1.338917364670584
$row7[0,40] = 15.48282828282828
$row7[0,41] = -0.09751218902362796
$row7[0,42] = 15.48282828282828
$ws.Range("A7:AQ7").Value = $row7

# --- Row 8 ---
$row8 = New-Object 'object[,]' 1,43
$row8[0,0] = "United Kingdom"
$row8[0,1] = "Beazley plc (LSE:BEZ)"
$row8[0,2] = "Insurance (Prop/Cas.)"
$row8[0,3] = 0.0789
$row8[0,4] = -0.189
$row8[0,5] = 1.65
$row8[0,6] = 0.1145946757781699
$row8[0,7] = 0.1145946757781699
$row8[0,8] = 0.04633349603874892
$row8[0,9] = 0.04384472539438184
$row8[0,10] = 82.8
$row8[0,11] = 0.03108925017835016
$row8[0,12] = 97.5
$row8[0,13] = 0.03239417901521695
$row8[0,14] = 1.177536231884058
$row8[0,15] = 77
$row8[0,16] = 0.0255830952222739
$row8[0,17] = 0.9299516908212561
$row8[0,18] = 20.5
$row8[0,19] = 0.2102564102564103
$row8[0,20] = 350.3
$row8[0,21] = 0.1163864708618513
$row8[0,22] = 0.0530871321407963
$row8[0,23] = 0.04974242011147174
$row8[0,24] = 0.003344712029324556
$row8[0,25] = 1.593263938741326
$row8[0,26] = 0.06985621987488463
$row8[0,27] = 0.04470574091029064
$row8[0,28] = 0.02515047896459399
$row8[0,29] = 610.6
$row8[0,30] = 0
$row8[0,31] = 610.6
$row8[0,32] = 260.3
$row8[0,33] = 0.1686553971936803
$row8[0,34] = 0.2504306455581987
$row8[0,35] = 0.07960001223204183
$row8[0,36] = 0.1246707217778629
$row8[0,37] = 35.2
$row8[0,38] = 35.2
$row8[0,39] = 4.460189919649379
$row8[0,40] = 3.505681818181818
$row8[0,41] = 1.901387874360847
$row8[0,42] = 3.505681818181818
$ws.Range("A8:AQ8").Value = $row8

